$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1")

$boldPart = "Table 1"
$restPart = " Characteristics of coastal segments used to evaluate seagrass depth of colonization estimates (see Fig. 2 for spatial distribution).  Year is the date of the seagrass coverage and bathymetric data.  Latitude and longitude are the geographic centers of each segment.  Area and depth values are square kilometers and meters, respectively.  Secchi measurements (m) were obtained from the Florida Department of Environmental Protection's Impaired Waters Rule (IWR) database, update number 40.  Secchi mean and standard errors are based on all observations within the ten years preceding each seagrass survey"
$fullText = $boldPart + $restPart

$rng = $ws.Range("A1")
$rng.Value = $fullText

# First run: "Table 1" in bold, Times 12pt (matches the rest of the caption's font).
$boldChars = $rng.Characters(1, $boldPart.Length)
$boldChars.Font.Name = "Times"
$boldChars.Font.Size = 12
$boldChars.Font.Bold = $true
$boldChars.Font.ColorIndex = -4105

# Second run: the remainder of the caption, regular weight, same font/size.
$restChars = $rng.Characters($boldPart.Length + 1, $restPart.Length)
$restChars.Font.Name = "Times"
$restChars.Font.Size = 12
$restChars.Font.Bold = $false
$restChars.Font.ColorIndex = -4105

Write-Output "done"
